# Apply the "Updated model policy schedule" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("About")

$ws1.Cells.Item(11,1).Value2 = "As of EPS 3.1.0, this lever supports the three energy carriers (electricity,"
$ws1.Cells.Item(12,1).Value2 = "district heat, and hydrogen), as well as fuels produced by the natural gas"
$ws1.Cells.Item(13,1).Value2 = "and petroleum, coal, biomass, and biofuel industries, as noted on the blue tab."
$ws1.Cells.Item(14,1).Value2 = ""
$ws1.Cells.Item(15,1).Value2 = "In the U.S. model, by default, we allow the suppliers of energy carriers"
$ws1.Cells.Item(16,1).Value2 = "(electricity, district heat, and hydrogen) to pass through changes in their"
$ws1.Cells.Item(17,1).Value2 = "expenses, while other fuel suppliers do not, due to the influence of a global"
$ws1.Cells.Item(18,1).Value2 = "market on setting prices."

# ---------------------------------------------------------------------------
# Sheet "BAEPAbCiPC"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("BAEPAbCiPC")

$ws2.Cells.Item(1,1).Value2 = "Unit: boolean (1 or 0)"
$ws2.Cells.Item(1,1).Font.Italic = $true

$ws2.Cells.Item(5,1).Value2 = "nuclear (NOT USED)"
$ws2.Range("A5:B5").Interior.ColorIndex = -4142

$ws2.Cells.Item(21,1).Value2 = "municipal solid waste (NOT USED)"
$ws2.Range("A21:B21").Interior.ColorIndex = -4142

$ws2.Range("A2:B2").Interior.ColorIndex = -4142
$ws2.Range("A15:B15").Interior.ColorIndex = -4142
$ws2.Range("A22:B22").Interior.ColorIndex = -4142

$ws2.PageSetup.Orientation = 1
